$wb = $excel.ActiveWorkbook

$wsScopes   = $wb.Worksheets.Item("#Scopes")
$wsConcerns = $wb.Worksheets.Item("#Concerns")

# --- Rename "Concerns" related shared strings to "Objectives" ------------
# Order matters for where new unique shared-string entries land in the table:
# we first touch B1 (ttIsaConcern -> ttIsaObjective) then A1 ([Concerns] -> [Objectives])
# so that the renamed string keeps the earlier slot and the brand new one is appended last.
$wsConcerns.Cells.Item(1, 2).Value = "ttIsaObjective"
$wsConcerns.Cells.Item(1, 1).Value = "[Objectives]"

# --- Update the Conc_ -> Obj_ concatenation formulas in column A ---------
for ($r = 3; $r -le 21; $r++) {
    $formula = '=IF(OR($C' + $r + '="",$D' + $r + '=""),"",CONCATENATE("Obj_",$C' + $r + ',"_",$D' + $r + '))'
    $wsConcerns.Cells.Item($r, 1).Formula = $formula
}

# --- Hide column B on the #Concerns sheet ---------------------------------
$wsConcerns.Columns.Item(2).Hidden = $true

# --- Update sheet selections / active sheet -------------------------------
$wsConcerns.Activate()
$wsConcerns.Range("D5").Select()
